# calculate player's power to basicinfo.power
# Adds a new "INT_power" column to every HouseFunction sheet (dwelling,
# woodcutter, quarrier, miner, farmer), reusing the style pair (header/data)
# that the "farmer" sheet's production column already uses (s=4 / s=3).

$wb = $excel.ActiveWorkbook

$powerValues = @(0, 2, 4, 12, 24, 72, 120, 168, 216, 264, 312, 360, 408, 456, 504, 552, 600, 648, 696, 744, 792)

# sheet index -> column letter that will receive the new "INT_power" data
$targets = @{
    1 = "D"   # dwelling: A,B,C already used -> new col D
    2 = "C"   # woodcutter: A,B already used -> new col C
    3 = "C"   # quarrier
    4 = "C"   # miner
    5 = "C"   # farmer
}

# Sheet "farmer" (index 5) column B already carries the header/data style
# pair (s="4" header, s="3" data) that the new power column should use.
$styleSourceSheet = $wb.Worksheets.Item(5)
$styleSourceRange = "B1:B22"

foreach ($idx in 1..5) {
    $ws = $wb.Worksheets.Item($idx)
    $col = $targets[$idx]

    # Copy the style template (header + 21 data rows) into the new column.
    $styleSourceSheet.Range($styleSourceRange).Copy()
    $ws.Range("$col" + "1:" + "$col" + "22").PasteSpecial(-4122)

    # Header text.
    $ws.Range("$col" + "1").Value = "INT_power"

    # Data values (rows 2..22).
    for ($i = 0; $i -lt $powerValues.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, $ws.Range("$col" + "1").Column).Value = $powerValues[$i]
    }
}

# Restore / set the selections on each sheet to match the edited ranges.
$wb.Worksheets.Item(1).Range("D1:D22").Select() | Out-Null
$wb.Worksheets.Item(2).Range("C1:C22").Select() | Out-Null
$wb.Worksheets.Item(3).Range("C1:C22").Select() | Out-Null
$wb.Worksheets.Item(4).Range("C1:C22").Select() | Out-Null
# Farmer stays the active/selected sheet, with a single-cell selection.
$wb.Worksheets.Item(5).Range("D19").Select() | Out-Null
